$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -3
